# Several more weeks of data
# Adds weekly rows for w/o 2020-08-10, 08-17, 08-24, 08-31 (Excel serials 44061,44068,44075,44081)
# to Fallecido_Recuperado and Provincias_Semanal, and moves the active sheet
# from Por_Edad to Provincias_Semanal.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Fallecido_Recuperado": 4 new weekly summary rows (22-25)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Fallecido_Recuperado")

$lastRow1 = 21
$newRows1 = 4

# Copy the number-format (date style) of the last existing row down into the
# new rows so column A keeps its date formatting.
$ws1.Range("A" + $lastRow1).Copy()
$ws1.Range("A22:A25").PasteSpecial(-4122)

$data1 = @(
    @(44061, 88127, 1501, 56760),
    @(44068, 92557, 1613, 63478),
    @(44075, 95627, 1765, 69519),
    @(44081, 100131, 1889, 73795)
)

for ($i = 0; $i -lt $data1.Length; $i++) {
    $r = $lastRow1 + 1 + $i
    $row = $data1[$i]
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
}

# ---------------------------------------------------------------------------
# Sheet "Provincias_Semanal": 4 new weeks x 32 provinces (rows 642-769)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Provincias_Semanal")

$lastRow2 = 641
$ws2.Range("A" + $lastRow2).Copy()
$ws2.Range("A642:A769").PasteSpecial(-4122)

$data2 = @(
    @(44061, "Distrito Nacional", 2406.09, 278),
    @(44061, "Azua", 709.99, 17),
    @(44061, "Baoruco", 327.71, 3),
    @(44061, "Barahona", 706.75, 15),
    @(44061, "Dajabon", 246.74, 5),
    @(44061, "Duarte", 679.8, 100),
    @(44061, "Elias Pina", 328.07, 2),
    @(44061, "El Seibo", 343.29, 2),
    @(44061, "Espaillat", 611.65, 42),
    @(44061, "Independencia", 332.06, 5),
    @(44061, "La Altagracia", 883.12, 22),
    @(44061, "La Romana", 692.96, 30),
    @(44061, "La Vega", 753.48, 68),
    @(44061, "Maria Trinidad Sanchez", 712.28, 7),
    @(44061, "Monte Cristi", 187.68, 6),
    @(44061, "Pedernales", 1231.53, 3),
    @(44061, "Peravia", 445.72, 26),
    @(44061, "Puerto Plata", 469.06, 82),
    @(44061, "Hermanas Mirabal", 417.81, 22),
    @(44061, "Samana", 214.87, 3),
    @(44061, "San Cristobal", 527.9, 82),
    @(44061, "San Juan", 555.16, 17),
    @(44061, "San Pedro de Macoris", 410.87, 22),
    @(44061, "Sanchez Ramirez", 800.12, 14),
    @(44061, "Santiago", 886.36, 194),
    @(44061, "Santiago Rodriguez", 640.87, 7),
    @(44061, "Valverde", 232.57, 22),
    @(44061, "Monsenor Nouel", 530.76, 17),
    @(44061, "Monte Plata", 104.17, 12),
    @(44061, "Hato Mayor", 171.4, 7),
    @(44061, "San Jose de Ocoa", 595.53, 11),
    @(44061, "Santo Domingo", 661.53, 358),
    @(44068, "Distrito Nacional", 2531.67, 304),
    @(44068, "Azua", 732.49, 19),
    @(44068, "Baoruco", 392.06, 3),
    @(44068, "Barahona", 723.13, 16),
    @(44068, "Dajabon", 246.74, 5),
    @(44068, "Duarte", 732.65, 104),
    @(44068, "Elias Pina", 329.64, 2),
    @(44068, "El Seibo", 347.65, 2),
    @(44068, "Espaillat", 620.4, 49),
    @(44068, "Independencia", 333.77, 5),
    @(44068, "La Altagracia", 931.51, 23),
    @(44068, "La Romana", 736.99, 31),
    @(44068, "La Vega", 794.82, 72),
    @(44068, "Maria Trinidad Sanchez", 747, 8),
    @(44068, "Monte Cristi", 223.51, 8),
    @(44068, "Pedernales", 1231.53, 3),
    @(44068, "Peravia", 489.78, 30),
    @(44068, "Puerto Plata", 496.97, 85),
    @(44068, "Hermanas Mirabal", 453.62, 22),
    @(44068, "Samana", 232.7, 3),
    @(44068, "San Cristobal", 552.38, 89),
    @(44068, "San Juan", 608.83000000000004, 19),
    @(44068, "San Pedro de Macoris", 432.18, 24),
    @(44068, "Sanchez Ramirez", 866.57, 15),
    @(44068, "Santiago", 946.35, 211),
    @(44068, "Santiago Rodriguez", 658.33, 7),
    @(44068, "Valverde", 238.23, 24),
    @(44068, "Monsenor Nouel", 556.01, 17),
    @(44068, "Monte Plata", 108.88, 13),
    @(44068, "Hato Mayor", 185.4, 7),
    @(44068, "San Jose de Ocoa", 663.54, 11),
    @(44068, "Santo Domingo", 685.2, 382),
    @(44075, "Distrito Nacional", 2609.41, 328),
    @(44075, "Azua", 754.54, 20),
    @(44075, "Baoruco", 395.03, 4),
    @(44075, "Barahona", 725.25, 16),
    @(44075, "Dajabon", 249.75, 5),
    @(44075, "Duarte", 781.83, 106),
    @(44075, "Elias Pina", 332.8, 4),
    @(44075, "El Seibo", 348.63, 2),
    @(44075, "Espaillat", 632.91, 56),
    @(44075, "Independencia", 333.77, 5),
    @(44075, "La Altagracia", 960.09, 31),
    @(44075, "La Romana", 750.56, 33),
    @(44075, "La Vega", 839.55, 78),
    @(44075, "Maria Trinidad Sanchez", 759.05, 9),
    @(44075, "Monte Cristi", 225.22, 8),
    @(44075, "Pedernales", 1234.3900000000001, 3),
    @(44075, "Peravia", 507, 36),
    @(44075, "Puerto Plata", 523.08000000000004, 90),
    @(44075, "Hermanas Mirabal", 472.07, 22),
    @(44075, "Samana", 236.27, 3),
    @(44075, "San Cristobal", 564.14, 99),
    @(44075, "San Juan", 641.75, 24),
    @(44075, "San Pedro de Macoris", 456.77, 31),
    @(44075, "Sanchez Ramirez", 910.66, 16),
    @(44075, "Santiago", 1003.47, 242),
    @(44075, "Santiago Rodriguez", 660.76, 8),
    @(44075, "Valverde", 238.8, 24),
    @(44075, "Monsenor Nouel", 600.76, 18),
    @(44075, "Monte Plata", 113.07, 13),
    @(44075, "Hato Mayor", 188.89, 9),
    @(44075, "San Jose de Ocoa", 667.22, 11),
    @(44075, "Santo Domingo", 703.13, 411),
    @(44081, "Distrito Nacional", 2684.95, 343),
    @(44081, "Azua", 773.43, 22),
    @(44081, "Baoruco", 398, 5),
    @(44081, "Barahona", 751.15, 17),
    @(44081, "Dajabon", 249.75, 5),
    @(44081, "Duarte", 830.01, 109),
    @(44081, "Elias Pina", 334.37, 4),
    @(44081, "El Seibo", 368.95, 4),
    @(44081, "Espaillat", 683.77, 61),
    @(44081, "Independencia", 383.4, 5),
    @(44081, "La Altagracia", 1015.55, 32),
    @(44081, "La Romana", 780.27, 36),
    @(44081, "La Vega", 888.42, 80),
    @(44081, "Maria Trinidad Sanchez", 771.1, 9),
    @(44081, "Monte Cristi", 276.39999999999998, 9),
    @(44081, "Pedernales", 1240.1099999999999, 3),
    @(44081, "Peravia", 534.36, 39),
    @(44081, "Puerto Plata", 554.59, 95),
    @(44081, "Hermanas Mirabal", 490.52, 23),
    @(44081, "Samana", 254.99, 3),
    @(44081, "San Cristobal", 587.20000000000005, 105),
    @(44081, "San Juan", 699.03, 32),
    @(44081, "San Pedro de Macoris", 485.96, 44),
    @(44081, "Sanchez Ramirez", 948.82, 16),
    @(44081, "Santiago", 1081.26, 273),
    @(44081, "Santiago Rodriguez", 665.32, 8),
    @(44081, "Valverde", 279.54000000000002, 26),
    @(44081, "Monsenor Nouel", 639.78, 19),
    @(44081, "Monte Plata", 121.97, 17),
    @(44081, "Hato Mayor", 204.05, 10),
    @(44081, "San Jose de Ocoa", 707.66, 12),
    @(44081, "Santo Domingo", 723.43, 423)
)

for ($i = 0; $i -lt $data2.Length; $i++) {
    $r = $lastRow2 + 1 + $i
    $row = $data2[$i]
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
}

# ---------------------------------------------------------------------------
# Selections / active sheet, matching the saved workbook view state
# ---------------------------------------------------------------------------
$ws1.Range("D26").Select()

# Provincias_Semanal becomes the active sheet/tab (was Por_Edad before).
$ws2.Activate()
$ws2.Range("D770").Select()
